$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) so numeric-looking strings (e.g. "1.001")
# are written back as literal text instead of being coerced into numbers,
# matching the inline-string storage used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.487.90'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.849.16'
$ws.Range("E3").Value = '  -1.34%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '260.46'
$ws.Range("E5").Value = '  -7.98%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.5161'
$ws.Range("E7").Value = '  -0.59%  '
$ws.Range("D8").Value = '0.3248'
$ws.Range("E8").Value = '  -7.96%  '
$ws.Range("D9").Value = '0.06765'
$ws.Range("E9").Value = '  -4.73%  '
$ws.Range("D10").Value = '18.89'
$ws.Range("E10").Value = '  -6.64%  '
$ws.Range("D11").Value = '0.7712'
$ws.Range("E11").Value = '  -6.09%  '
$ws.Range("D12").Value = '0.07718'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '1.889.61'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '88.60'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("E15").Value = '  -2.76%  '
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("E17").Value = '  -2.49%  '
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '0.000007920'
$ws.Range("E19").Value = '  -2.98%  '
$ws.Range("D20").Value = '26.537.41'
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").Value = '2.096.38'
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").Value = '4.522'
$ws.Range("D23").Value = '9.531'
$ws.Range("E23").Value = '  -6.39%  '
$ws.Range("D24").Value = '5.918'
$ws.Range("E24").Value = '  -5.17%  '
$ws.Range("D25").Value = '2.350'
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").Value = '144.50'
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("E28").Value = '  -2.65%  '
$ws.Range("D29").Value = '111.21'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").Value = '4.208'
$ws.Range("E30").Value = '  -4.86%  '
$ws.Range("D31").Value = '4.174'
$ws.Range("E31").Value = '  -4.32%  '
$ws.Range("D32").Value = '0.08758'
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = '0.04813'
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").Value = '1.133'
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("D35").Value = '2.843'
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("D36").Value = '0.6868'
$ws.Range("E36").Value = '  -8.06%  '
$ws.Range("D37").Value = '3.118'
$ws.Range("E37").Value = '  -5.47%  '
$ws.Range("D38").Value = '0.01796'
$ws.Range("E38").Value = '  -4.57%  '
$ws.Range("D39").Value = '2.205'
$ws.Range("E39").Value = '  -8.81%  '
$ws.Range("D40").Value = '0.4891'
$ws.Range("E40").Value = '  -8.20%  '
$ws.Range("D41").Value = '113.03'
$ws.Range("E41").Value = '  -2.95%  '
$ws.Range("D42").Value = '0.8997'
$ws.Range("E42").Value = '  -7.73%  '
$ws.Range("D43").Value = '6.136'
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").Value = '0.9999'
$ws.Range("D45").Value = '7.779'
$ws.Range("E45").Value = '  -5.11%  '
$ws.Range("D46").Value = '0.4221'
$ws.Range("E46").Value = '  -8.72%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.155'
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1257'
$ws.Range("E48").Value = '  -8.25%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05893'
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").Value = '35.20'
$ws.Range("E50").Value = '  -4.13%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '59.25'
$ws.Range("E51").Value = '  -4.35%  '

# Restore the default cell style so no stray formatting is introduced
# by the temporary text number-format above.
$ws.Range("D2:D51").Style = "Normal"

